$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new rows of unit metadata beneath the existing table
$ws.Range("A7").Value = "age"
$ws.Range("B7").Value = "Ar Age"
$ws.Range("C7").Value = "kyr"

$ws.Range("A8").Value = "concentration"
$ws.Range("B8").Value = "Particle Concentration "
$ws.Range("C8").Value = "ppb"

$ws.Range("A9").Value = "tac"
$ws.Range("B9").Value = "Total Air Content"
$ws.Range("C9").Value = "cm`$^3`$/g"

$ws.Range("F15").Select()
